$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 186, shifting existing rows 186:310 down to 187:311
$ws.Rows.Item(186).Insert()

# Populate the new row 186 with the new data record
$ws.Range("A186").Value = 7
$ws.Range("B186").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C186").Value = "Ñuble"
$ws.Range("D186").Value = 44777
$ws.Range("E186").Value = 16
$ws.Range("F186").Value = 100114013
$ws.Range("G186").Value = "Zanahoria"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 120
$ws.Range("K186").Value = 9000
$ws.Range("L186").Value = 10000
$ws.Range("M186").Value = 9500
$ws.Range("N186").Value = "`$/saco 20 kilos"
$ws.Range("O186").Value = "Provincia de Diguillín"
$ws.Range("P186").Value = 475
$ws.Range("Q186").Value = 20
$ws.Range("R186").Value = "Hortaliza"
